$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Pad the sheet with formatted (but empty) rows 6-17 FIRST, while column B
#    still carries its original "vertical top" only style (B3/B4) - this is
#    the style the blank filler rows should inherit.
# ---------------------------------------------------------------------------
$ws.Range("B3").Copy()
$ws.Range("A6:C17").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2. Add the brand-new last row ("Register workspace"), then reword the
#    "Register data" row's description.
# ---------------------------------------------------------------------------

# -- New row 5 ("Register workspace") ----------------------------------------
$ws.Range("A5").Value = 44183
$ws.Range("B5").Value = "Register workspace"

# -- Row 2 ("Register data") -------------------------------------------------
$ws.Range("C2").Value = "System has to store the data time that user enters through the GUI"

# -- finish row 5's description -----------------------------------------------
$ws.Range("C5").Value = "System has to allow that de user can create a workspace in which its store his personal configuration for that each time  that user enter to the program its load automatically that workspace"

# ---------------------------------------------------------------------------
# 3. Formatting clean-up: make every "Specification date" cell share the same
#    date format + top alignment, every "Name" cell left/top aligned, and
#    every "Description" cell top aligned + wrapped.
# ---------------------------------------------------------------------------

# Column A (dates): copy the existing date style (numFmt 14 + vertical top)
# from A3 onto A2 and the new A5 so no new number format gets created.
$ws.Range("A3").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$ws.Range("A5").PasteSpecial(-4122)

# Column C (descriptions): copy the existing top+wrap style from C4 onto
# C2, C3 and the new C5.
$ws.Range("C4").Copy()
$ws.Range("C2").PasteSpecial(-4122)
$ws.Range("C3").PasteSpecial(-4122)
$ws.Range("C5").PasteSpecial(-4122)

# Column B (names): apply left + top alignment to every data row (vertical
# first so the intermediate state reuses the already-existing "top" style
# instead of minting a throw-away "left only" one).
$ws.Range("B2:B5").VerticalAlignment = -4160
$ws.Range("B2:B5").HorizontalAlignment = -4131

# Row heights for the taller, wrapped rows.
$ws.Rows.Item(2).RowHeight = 30
$ws.Rows.Item(5).RowHeight = 60

# ---------------------------------------------------------------------------
# 4. Restore the selection to the cell the author last had selected.
# ---------------------------------------------------------------------------
$ws.Range("C9").Select()

Write-Host "done"
